# "flappy bird part I"
# Adds three new notes to the "Emelt" (main) tanmenet sheet:
#   C5 = "Flappy bird clouds"
#   E6 = "Hangosítás 17:00-től"
#   E9 = "Tantestületi ebéd 14:00-től"
# Updates the zoom level / selected cell on both sheets, and normalises the
# header/footer font-style label from "Regular" to "Általános".

$wb = $excel.ActiveWorkbook

$wsEmelt = $wb.Worksheets.Item("Emelt")
$wsKozep = $wb.Worksheets.Item("Közép")

# --- new content on the main sheet ---------------------------------------
$wsEmelt.Range("C5").Value = "Flappy bird clouds"
$wsEmelt.Range("E6").Value = "Hangosítás 17:00-től"
$wsEmelt.Range("E9").Value = "Tantestületi ebéd 14:00-től"

# --- header / footer text on both sheets ----------------------------------
$wsEmelt.PageSetup.CenterHeader = '&"Times New Roman,Általános"&12&A'
$wsEmelt.PageSetup.CenterFooter = '&"Times New Roman,Általános"&12Page &P'

$wsKozep.PageSetup.CenterHeader = '&"Times New Roman,Általános"&12&A'
$wsKozep.PageSetup.CenterFooter = '&"Times New Roman,Általános"&12Page &P'

# --- view state: zoom + selected cell on the "Közép" sheet ----------------
$wsKozep.Activate()
$excel.ActiveWindow.Zoom = 120
$wsKozep.Range("D3").Select()

# --- view state: zoom + selected cell on the "Emelt" sheet (left active) --
$wsEmelt.Activate()
$excel.ActiveWindow.Zoom = 120
$wsEmelt.Range("C6").Select()
